$wb = $excel.ActiveWorkbook

# --- Chart sheet: fix D32 (Impressions for 2025-11-04) to be numeric 25 ---
$wsChart = $wb.Worksheets.Item("Chart")
$wsChart.Range("D32").Value = 25

# --- Critical issues sheet: reorder data rows 5-8 ---
$wsCritical = $wb.Worksheets.Item("Critical issues")

$wsCritical.Range("A5").Value = "Excluded by ‘noindex’ tag"
$wsCritical.Range("B5").Value = "Website"
$wsCritical.Range("C5").Value = "Not Started"
$wsCritical.Range("D5").Value = 13

$wsCritical.Range("A6").Value = "Duplicate, Google chose different canonical than user"
$wsCritical.Range("B6").Value = "Google systems"
$wsCritical.Range("C6").Value = "Started"
$wsCritical.Range("D6").Value = 42

$wsCritical.Range("A7").Value = "Discovered - currently not indexed"
$wsCritical.Range("B7").Value = "Google systems"
$wsCritical.Range("C7").Value = "Started"
$wsCritical.Range("D7").Value = 6

$wsCritical.Range("A8").Value = "Crawled - currently not indexed"
$wsCritical.Range("B8").Value = "Google systems"
$wsCritical.Range("C8").Value = "Started"
$wsCritical.Range("D8").Value = 4
